# Update scripts with new TPM values (NATMI LR-pairs output refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: recomputed specificities / weights ---
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3643233333333333
$ws.Range("N2").Value = 1.09297
$ws.Range("O2").Value = 0.09146135066097912
$ws.Range("P2").Value = 0.09146135066097912
$ws.Range("Q2").Value = 0.1740056816444445
$ws.Range("R2").Value = 1.5660511348
$ws.Range("S2").Value = 0.09146135066097912
$ws.Range("T2").Value = 0.09146135066097912

# --- Row 3: recomputed receptor specificities ---
$ws.Range("O3").Value = 0.2211540931751474
$ws.Range("P3").Value = 0.2211540931751474
$ws.Range("S3").Value = 0.2211540931751474
$ws.Range("T3").Value = 0.2211540931751474

# --- Row 4: recomputed receptor expression / specificities ---
$ws.Range("M4").Value = 2.721212
$ws.Range("N4").Value = 8.163636
$ws.Range("O4").Value = 0.6831451685449673
$ws.Range("P4").Value = 0.6831451685449673
$ws.Range("Q4").Value = 1.299687134026667
$ws.Range("R4").Value = 11.69718420624
$ws.Range("S4").Value = 0.6831451685449673
$ws.Range("T4").Value = 0.6831451685449673

# --- Row 5: new target cluster "Resolving-Mac" ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt2"
$ws.Range("C5").Value = "Fzd3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4776133333333334
$ws.Range("H5").Value = 1.43284
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.016887
$ws.Range("N5").Value = 0.050661
$ws.Range("O5").Value = 0.004239387618906157
$ws.Range("P5").Value = 0.004239387618906157
$ws.Range("Q5").Value = 0.008065456360000001
$ws.Range("R5").Value = 0.07258910724000001
$ws.Range("S5").Value = 0.004239387618906157
$ws.Range("T5").Value = 0.004239387618906157
